# Atualiza planilhas com concursos mais recentes
# Adds the 5 newest Mega-Sena draws (concursos 2955-2959) to the bottom
# of the "MEGA SENA" sheet, right after the existing last row (411).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Concurso, Bola1..Bola6
$novosConcursos = @(
    @(2955,  9, 13, 21, 32, 33, 59),
    @(2956, 10, 18, 21, 24, 43, 47),
    @(2957, 19, 28, 36, 37, 48, 52),
    @(2958,  7,  9, 14, 35, 42, 49),
    @(2959, 18, 26, 35, 41, 44, 45)
)

$linhaInicial = 412

for ($i = 0; $i -lt $novosConcursos.Count; $i++) {
    $linha = $linhaInicial + $i
    $valores = $novosConcursos[$i]
    for ($col = 0; $col -lt $valores.Count; $col++) {
        $ws.Cells.Item($linha, $col + 1).Value = $valores[$col]
    }
}

$ultimaLinha = $linhaInicial + $novosConcursos.Count - 1

# Mirror the author's final selection: B412:G416 with B412 active.
[void]$ws.Range("B412:G$ultimaLinha").Select()
